$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($ws, $addr, $val)
    $ws.Range($addr).Value = "'" + $val
    $ws.Range($addr).Style = "Normal"
}

Set-TextValue $ws "D2" "332.64"
Set-TextValue $ws "E2" "0.54%"
Set-TextValue $ws "D3" "41.29"
Set-TextValue $ws "E3" "0.54%"
Set-TextValue $ws "D4" "5.696"
Set-TextValue $ws "E4" "-0.88%"
Set-TextValue $ws "D5" "0.08431"
Set-TextValue $ws "E5" "3.54%"
Set-TextValue $ws "D6" "8.829"
Set-TextValue $ws "E6" "1.00%"
Set-TextValue $ws "E7" "0.05%"
Set-TextValue $ws "D8" "1.988"
Set-TextValue $ws "E8" "-2.90%"
Set-TextValue $ws "D9" "2.932"
Set-TextValue $ws "E9" "-1.07%"
Set-TextValue $ws "D10" "0.9282"
Set-TextValue $ws "E10" "0.68%"
Set-TextValue $ws "E11" "0.96%"
Set-TextValue $ws "D12" "0.1969"
Set-TextValue $ws "E12" "0.81%"
Set-TextValue $ws "D13" "0.09365"
Set-TextValue $ws "E13" "0.36%"
Set-TextValue $ws "D14" "0.04016"
Set-TextValue $ws "E14" "9.80%"
Set-TextValue $ws "D15" "0.1064"
Set-TextValue $ws "E15" "0.81%"
Set-TextValue $ws "D16" "0.001317"
Set-TextValue $ws "E16" "1.45%"
Set-TextValue $ws "D17" "0.006101"
Set-TextValue $ws "E17" "-0.99%"
Set-TextValue $ws "D18" "3.434"
Set-TextValue $ws "E18" "1.55%"
Set-TextValue $ws "D19" "0.3510"
Set-TextValue $ws "E19" "0.75%"
Set-TextValue $ws "D20" "9.169"
Set-TextValue $ws "E20" "10.55%"
Set-TextValue $ws "E21" "-3.68%"
Set-TextValue $ws "D22" "0.2630"
Set-TextValue $ws "E22" "-0.71%"
Set-TextValue $ws "D23" "0.04417"
Set-TextValue $ws "E23" "-0.44%"
Set-TextValue $ws "D24" "0.001244"
Set-TextValue $ws "E24" "-1.31%"
Set-TextValue $ws "D25" "0.004383"
Set-TextValue $ws "E25" "0.64%"
Set-TextValue $ws "E26" "-3.87%"
Set-TextValue $ws "D27" "0.0003999"
Set-TextValue $ws "E27" "0.15%"
Set-TextValue $ws "D39" "0.02813"
Set-TextValue $ws "E39" "1.07%"
Set-TextValue $ws "E40" "0.48%"
Set-TextValue $ws "D41" "0.007914"
Set-TextValue $ws "E41" "4.09%"
Set-TextValue $ws "D42" "0.1441"
Set-TextValue $ws "E42" "1.15%"
Set-TextValue $ws "D43" "0.008983"
Set-TextValue $ws "E43" "-9.63%"
Set-TextValue $ws "D44" "0.002085"
Set-TextValue $ws "E44" "-1.55%"
Set-TextValue $ws "D45" "0.01039"
Set-TextValue $ws "E45" "-12.32%"
Set-TextValue $ws "D46" "0.00007172"
Set-TextValue $ws "E46" "6.11%"
Set-TextValue $ws "E47" "0.26%"
Set-TextValue $ws "D48" "0.003411"
Set-TextValue $ws "E48" "15.93%"
Set-TextValue $ws "D49" "0.002283"
Set-TextValue $ws "E49" "0.19%"
Set-TextValue $ws "D50" "0.00002105"
Set-TextValue $ws "E50" "0.26%"
Set-TextValue $ws "D51" "0.0002005"
Set-TextValue $ws "E51" "0.26%"

Write-Host "Applied updates to 72 cells"
